$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit adds two new weekly price records for "Berenjena" (Femacal de
# La Calera) into the existing, reverse-chronological data table:
#   - one brand-new record inserted at the top of the table (row 300)
#   - one brand-new record inserted further down (originally landing at
#     row 356), which naturally shifts every following record down by one
#     row (the table keeps growing from the bottom: rows 419-421 now hold
#     what used to be rows 417-419).
# Using Rows(...).Insert() lets Excel perform that cascading shift for us
# instead of having to rewrite ~120 rows by hand.

# --- Insert and populate the first new record at row 300 ---
$ws.Rows(300).Insert()

$ws.Cells.Item(300, 1).Value = 3
$ws.Cells.Item(300, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(300, 3).Value = "Coquimbo"
$ws.Cells.Item(300, 4).Value = 45006
$ws.Cells.Item(300, 5).Value = 5
$ws.Cells.Item(300, 6).Value = 100112001
$ws.Cells.Item(300, 7).Value = "Berenjena"
$ws.Cells.Item(300, 8).Value = "Sin especificar"
$ws.Cells.Item(300, 9).Value = "Primera"
$ws.Cells.Item(300, 10).Value = 40
$ws.Cells.Item(300, 11).Value = 8000
$ws.Cells.Item(300, 12).Value = 8000
$ws.Cells.Item(300, 13).Value = 8000
$ws.Cells.Item(300, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(300, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(300, 16).Value = 133
$ws.Cells.Item(300, 17).Value = 60
$ws.Cells.Item(300, 18).Value = "Hortaliza"

# --- Insert and populate the second new record, which lands at row 356
#     (after the row-300 insert above has already shifted everything
#     from the old row 356 onward down by one) ---
$ws.Rows(356).Insert()

$ws.Cells.Item(356, 1).Value = 3
$ws.Cells.Item(356, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(356, 3).Value = "Coquimbo"
$ws.Cells.Item(356, 4).Value = 45005
$ws.Cells.Item(356, 5).Value = 5
$ws.Cells.Item(356, 6).Value = 100112001
$ws.Cells.Item(356, 7).Value = "Berenjena"
$ws.Cells.Item(356, 8).Value = "Sin especificar"
$ws.Cells.Item(356, 9).Value = "Primera"
$ws.Cells.Item(356, 10).Value = 100
$ws.Cells.Item(356, 11).Value = 6500
$ws.Cells.Item(356, 12).Value = 7000
$ws.Cells.Item(356, 13).Value = 6725
$ws.Cells.Item(356, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(356, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(356, 16).Value = 112
$ws.Cells.Item(356, 17).Value = 60
$ws.Cells.Item(356, 18).Value = "Hortaliza"
